# Insert a new record row for "Feria Lagunitas de Puerto Montt - Tomate" weekly update.
# This shifts existing rows 478..511 down to 479..512 and fills the newly
# inserted row 478 with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 478 (pushes rows 478-511 down to 479-512).
$ws.Rows.Item(478).Insert()

# Populate the new row 478 with the new record.
$ws.Cells.Item(478, 1).Value2  = 4
$ws.Cells.Item(478, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(478, 3).Value2  = "Los Lagos"
$ws.Cells.Item(478, 4).Value2  = 44578
$ws.Cells.Item(478, 5).Value2  = 10
$ws.Cells.Item(478, 6).Value2  = 100112020
$ws.Cells.Item(478, 7).Value2  = "Tomate"
$ws.Cells.Item(478, 8).Value2  = "Larga vida"
$ws.Cells.Item(478, 9).Value2  = "Extra"
$ws.Cells.Item(478, 10).Value2 = 300
$ws.Cells.Item(478, 11).Value2 = 20000
$ws.Cells.Item(478, 12).Value2 = 20000
$ws.Cells.Item(478, 13).Value2 = 20000
$ws.Cells.Item(478, 14).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(478, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(478, 16).Value2 = 1111
$ws.Cells.Item(478, 17).Value2 = 18
$ws.Cells.Item(478, 18).Value2 = "Hortaliza"
